$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 4 new sign-up rows into the milk table (between the old row 17
#    and the old closing/border row 18), pushing the sugar table down by 4
#    rows (old 18->22 .. old 24->28).
# ---------------------------------------------------------------------------
$ws.Range("A18:D21").EntireRow.Insert()

# Re-apply the existing formatting patterns used by the table (re-using the
# workbook's existing cell styles rather than inventing new ones) by copying
# formats down from the rows directly above the insertion point.
$ws.Range("A14:D14").Copy() | Out-Null
$ws.Range("A18:D18").PasteSpecial(-4122) | Out-Null

$ws.Range("A15:D17").Copy() | Out-Null
$ws.Range("A19:D21").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Row heights: the milk sign-up rows (5-22, everything up to and
#    including the closing bordered row) shrink from 40 to 27, and the
#    sugar sign-up data rows (now 26-28) shrink from 40 to 35.
# ---------------------------------------------------------------------------
for ($r = 5; $r -le 22; $r++) {
    $ws.Rows.Item($r).RowHeight = 27
}
for ($r = 26; $r -le 28; $r++) {
    $ws.Rows.Item($r).RowHeight = 35
}

# ---------------------------------------------------------------------------
# 3. Column A widens slightly.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.5

# ---------------------------------------------------------------------------
# 4. View state: scroll so row 5 is at the top and select C19.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()

# ---------------------------------------------------------------------------
# 5. Page setup: margins and print copies.
# ---------------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.LeftMargin = 50.4
$ps.RightMargin = 50.4
$ps.TopMargin = 54
$ps.BottomMargin = 54
$ps.HeaderMargin = 21.6
$ps.FooterMargin = 21.6
$ps.Copies = 3
